$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Daily_Data")

# New daily rows for 2026-02-10 (serial 46063), appended after existing row 573
$arr = New-Object 'object[,]' 22,8
$arr[0,0] = 46063
$arr[0,1] = 'ASAHI DEPOSITORY LLC Registered'
$arr[0,2] = 0
$arr[0,3] = 0
$arr[0,4] = 0
$arr[0,5] = 0
$arr[0,6] = 0
$arr[0,7] = 0
$arr[1,0] = 46063
$arr[1,1] = 'ASAHI DEPOSITORY LLC Eligible'
$arr[1,2] = 0
$arr[1,3] = 0
$arr[1,4] = 0
$arr[1,5] = 0
$arr[1,6] = 0
$arr[1,7] = 0
$arr[2,0] = 46063
$arr[2,1] = 'BRINK''S, INC. Registered'
$arr[2,2] = 73354.783
$arr[2,3] = 0
$arr[2,4] = 0
$arr[2,5] = 0
$arr[2,6] = 0
$arr[2,7] = 73354.783
$arr[3,0] = 46063
$arr[3,1] = 'BRINK''S, INC. Eligible'
$arr[3,2] = 85821.84699999999
$arr[3,3] = 0
$arr[3,4] = 0
$arr[3,5] = 0
$arr[3,6] = 0
$arr[3,7] = 85821.84699999999
$arr[4,0] = 46063
$arr[4,1] = 'CNT DEPOSITORY, INC. Registered'
$arr[4,2] = 1246.06
$arr[4,3] = 0
$arr[4,4] = 0
$arr[4,5] = 0
$arr[4,6] = 0
$arr[4,7] = 1246.06
$arr[5,0] = 46063
$arr[5,1] = 'CNT DEPOSITORY, INC. Eligible'
$arr[5,2] = 0
$arr[5,3] = 0
$arr[5,4] = 0
$arr[5,5] = 0
$arr[5,6] = 0
$arr[5,7] = 0
$arr[6,0] = 46063
$arr[6,1] = 'DELAWARE DEPOSITORY Registered'
$arr[6,2] = 1633.941
$arr[6,3] = 0
$arr[6,4] = 0
$arr[6,5] = 0
$arr[6,6] = 0
$arr[6,7] = 1633.941
$arr[7,0] = 46063
$arr[7,1] = 'DELAWARE DEPOSITORY Eligible'
$arr[7,2] = 18459.584
$arr[7,3] = 0
$arr[7,4] = 0
$arr[7,5] = 0
$arr[7,6] = 0
$arr[7,7] = 18459.584
$arr[8,0] = 46063
$arr[8,1] = 'HSBC BANK, USA Registered'
$arr[8,2] = 1394.758
$arr[8,3] = 0
$arr[8,4] = 0
$arr[8,5] = 0
$arr[8,6] = 0
$arr[8,7] = 1394.758
$arr[9,0] = 46063
$arr[9,1] = 'HSBC BANK, USA Eligible'
$arr[9,2] = 9281.978999999999
$arr[9,3] = 0
$arr[9,4] = 0
$arr[9,5] = 0
$arr[9,6] = 0
$arr[9,7] = 9281.978999999999
$arr[10,0] = 46063
$arr[10,1] = 'INTERNATIONAL DEPOSITORY SERVICES OF DELAWARE Registered'
$arr[10,2] = 2395.448
$arr[10,3] = 0
$arr[10,4] = 0
$arr[10,5] = 0
$arr[10,6] = 0
$arr[10,7] = 2395.448
$arr[11,0] = 46063
$arr[11,1] = 'INTERNATIONAL DEPOSITORY SERVICES OF DELAWARE Eligible'
$arr[11,2] = 0
$arr[11,3] = 0
$arr[11,4] = 0
$arr[11,5] = 0
$arr[11,6] = 0
$arr[11,7] = 0
$arr[12,0] = 46063
$arr[12,1] = 'JP MORGAN CHASE BANK NA Registered'
$arr[12,2] = 114061.421
$arr[12,3] = 0
$arr[12,4] = 0
$arr[12,5] = 0
$arr[12,6] = 0
$arr[12,7] = 114061.421
$arr[13,0] = 46063
$arr[13,1] = 'JP MORGAN CHASE BANK NA Eligible'
$arr[13,2] = 76408.66899999999
$arr[13,3] = 0
$arr[13,4] = 0
$arr[13,5] = 0
$arr[13,6] = 0
$arr[13,7] = 76408.66899999999
$arr[14,0] = 46063
$arr[14,1] = 'LOOMIS INTERNATIONAL (US) LLC Registered'
$arr[14,2] = 61157.444
$arr[14,3] = 0
$arr[14,4] = 0
$arr[14,5] = 0
$arr[14,6] = 0
$arr[14,7] = 61157.444
$arr[15,0] = 46063
$arr[15,1] = 'LOOMIS INTERNATIONAL (US) LLC Eligible'
$arr[15,2] = 71594.18700000001
$arr[15,3] = 0
$arr[15,4] = 0
$arr[15,5] = 0
$arr[15,6] = 0
$arr[15,7] = 71594.18700000001
$arr[16,0] = 46063
$arr[16,1] = 'MALCA-AMIT USA, LLC Registered'
$arr[16,2] = 395.145
$arr[16,3] = 0
$arr[16,4] = 0
$arr[16,5] = 0
$arr[16,6] = 0
$arr[16,7] = 395.145
$arr[17,0] = 46063
$arr[17,1] = 'MALCA-AMIT USA, LLC Eligible'
$arr[17,2] = 0
$arr[17,3] = 0
$arr[17,4] = 0
$arr[17,5] = 0
$arr[17,6] = 0
$arr[17,7] = 0
$arr[18,0] = 46063
$arr[18,1] = 'MANFRA, TORDELLA & BROOKES, LLC Registered'
$arr[18,2] = 49920.248
$arr[18,3] = 0
$arr[18,4] = 0
$arr[18,5] = 0
$arr[18,6] = 0
$arr[18,7] = 49920.248
$arr[19,0] = 46063
$arr[19,1] = 'MANFRA, TORDELLA & BROOKES, LLC Eligible'
$arr[19,2] = 2104.855
$arr[19,3] = 0
$arr[19,4] = 0
$arr[19,5] = 0
$arr[19,6] = 0
$arr[19,7] = 2104.855
$arr[20,0] = 46063
$arr[20,1] = 'STONEX PRECIOUS METALS LLC Registered'
$arr[20,2] = 14122.765
$arr[20,3] = 0
$arr[20,4] = 0
$arr[20,5] = 0
$arr[20,6] = 0
$arr[20,7] = 14122.765
$arr[21,0] = 46063
$arr[21,1] = 'STONEX PRECIOUS METALS LLC Eligible'
$arr[21,2] = 16.075
$arr[21,3] = 0
$arr[21,4] = 0
$arr[21,5] = 0
$arr[21,6] = 0
$arr[21,7] = 16.075

$rng = $ws.Range("A574:H595")
$rng.Value = $arr

# Match date formatting/style used by column A in the rest of the table (row 573)
$ws.Range("A574:A595").NumberFormat = $ws.Range("A573").NumberFormat

"Added rows 574-595"
